$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5, duplicating the current row 4 (which holds the
# older "Provincia de Diguillin" / $/atado record dated 2021-01-25) down to row 5,
# so that row 4 is freed up to hold the new weekly record dated 2022-05-11.
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(5).Insert()

# Update row 4 with the new weekly price entry.
$ws.Cells.Item(4, 4).Value() = 44692
$ws.Cells.Item(4, 10).Value() = 120
$ws.Cells.Item(4, 11).Value() = 3000
$ws.Cells.Item(4, 12).Value() = 3500
$ws.Cells.Item(4, 13).Value() = 3250
$ws.Cells.Item(4, 14).Value() = "$/docena de matas"
$ws.Cells.Item(4, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(4, 16).Value() = 542
$ws.Cells.Item(4, 17).Value() = 6
$ws.Cells.Item(4, 18).Value() = "Hortaliza"
